# "Generate Report for Handback"
#
# The localization handback finished for both target files
# (81ea6839-...md and c88ab46c-...md) in both locales (zh-cn, de-de).
# This updates the status, records the handback target/handoff file names
# + timestamps, adds hyperlinks to the newly-available "Latest Target File"
# entries, and widens the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- github blob URLs backing the existing A2/A3 hyperlinks on each sheet ---
$urlFile1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d04a7fc41c9d2d36301a99b7a550490764c3619/e2e/81ea6839-1b41-47fb-8dda-449cc13d760d.md"
$urlFile2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d04a7fc41c9d2d36301a99b7a550490764c3619/e2e/c88ab46c-1dfb-4bd4-8f9a-ded54b5fdb6d.md"

# ------------------------------------------------------------------
# 1) Overview sheet: status columns (zh-cn / de-de) + wider columns
# ------------------------------------------------------------------
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

# column E & F: ~30 chars wide now that they show the longer status text
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666668

# ------------------------------------------------------------------
# 2) zh-cn sheet: status, handback file/time, target-file hyperlinks
# ------------------------------------------------------------------
$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus

# Latest Target File (I) -> hyperlink to the source .md, same address the
# A column already links to for that row
$ws2.Hyperlinks.Add($ws2.Range("I2"), $urlFile1, "", "", "81ea6839-1b41-47fb-8dda-449cc13d760d.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), $urlFile2, "", "", "c88ab46c-1dfb-4bd4-8f9a-ded54b5fdb6d.md")

# Latest Handback File (J)
$ws2.Range("J2").Value = "81ea6839-1b41-47fb-8dda-449cc13d760d.f7b20414592f92bb6db8669637f9e7531747813d.zh-cn.xlf"
$ws2.Range("J3").Value = "c88ab46c-1dfb-4bd4-8f9a-ded54b5fdb6d.169e821a57bd0f499f5be1bb4ce7805705c0a1e1.zh-cn.xlf"

# Latest Handback DateTime (K)
$ws2.Range("K2").Value = "2016-09-05 10:32:26"
$ws2.Range("K3").Value = "2016-09-05 10:32:26"

# column widths: C widens for the longer status text, I/J widen for filenames
$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws2.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws2.Columns.Item(10).ColumnWidth = 39.166666666666664

# ------------------------------------------------------------------
# 3) de-de sheet: status, handback file/time, target-file hyperlinks
# ------------------------------------------------------------------
$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

$ws3.Hyperlinks.Add($ws3.Range("I2"), $urlFile1, "", "", "81ea6839-1b41-47fb-8dda-449cc13d760d.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), $urlFile2, "", "", "c88ab46c-1dfb-4bd4-8f9a-ded54b5fdb6d.md")

$ws3.Range("J2").Value = "81ea6839-1b41-47fb-8dda-449cc13d760d.f7b20414592f92bb6db8669637f9e7531747813d.de-de.xlf"
$ws3.Range("J3").Value = "c88ab46c-1dfb-4bd4-8f9a-ded54b5fdb6d.169e821a57bd0f499f5be1bb4ce7805705c0a1e1.de-de.xlf"

$ws3.Range("K2").Value = "2016-09-05 10:32:35"
$ws3.Range("K3").Value = "2016-09-05 10:32:35"

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws3.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws3.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "Handback report generated."
